# The worksheet is protected, so it must be unprotected before any cell
# values can be changed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Unprotect()

# Update the confidential footer note: the "as of" date moves from
# 2021-04-28 to 2021-04-29 (row 7, column A).
$ws.Range("A7").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-29 for illustrative purposes only and are subject to change."

# Update the refreshed Weight / Percent Change figures.
$ws.Range("D2").Value = 0.8433413682032467
$ws.Range("E2").Value = 0.0003794106487922111

$ws.Range("D3").Value = 0.1566586317967534
$ws.Range("E3").Value = -0.00236020334059539

$ws.Range("E4").Value = -0.00004977353043644772

# Restore the sheet protection that was in place before the edits.
$ws.Protect()
